$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a blank column at E; it inherits column D formatting so the
# <cols> width group (D:K) stays merged after the edit.
$ws.Columns.Item(5).Insert()

# Shift the old column D (values+formats) into the new column E,
# which in turn pushed the old E:K along into F:L.
$ws.Range("D5:D102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# Populate column D with the newly reported period (the most recent column)
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1490800
$ws.Range("D9").Value = 819400
$ws.Range("D10").Value = 671400
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 91400
$ws.Range("D15").Value = 211400
$ws.Range("D17").Value = 1464900
$ws.Range("D18").Value = 25900
$ws.Range("D20").Value = -274800
$ws.Range("D21").Value = -37500
$ws.Range("D22").Value = "NA"
$ws.Range("D23").Value = -248900
$ws.Range("D24").Value = -5500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -243400
$ws.Range("D27").Value = -243400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 274800
$ws.Range("D33").Value = -243400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -243400
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 55300
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 177000
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 46700
$ws.Range("D46").Value = 279000
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 1870400
$ws.Range("D49").Value = 2290400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 97800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 4537600
$ws.Range("D57").Value = 89200
$ws.Range("D58").Value = 46600
$ws.Range("D59").Value = 311000
$ws.Range("D60").Value = 446800
$ws.Range("D61").Value = 3186700
$ws.Range("D62").Value = 489400
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 4122900
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -368300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 414700
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -243400
$ws.Range("D83").Value = 211400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 82400
$ws.Range("D91").Value = -77800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -2426300
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 2304300
$ws.Range("D101").Value = -6300
$ws.Range("D102").Value = -45900

# A handful of the shifted figures were also restated; apply those corrections
$ws.Range("E43").Value = 91500
$ws.Range("E45").Value = 34200
$ws.Range("E46").Value = 215300
$ws.Range("E54").Value = 1785400
$ws.Range("E59").Value = 157900
$ws.Range("E60").Value = 174900
$ws.Range("J62").Value = 10000
$ws.Range("E66").Value = 1554000
$ws.Range("F89").Value = 60400
$ws.Range("E91").Value = -56900
$ws.Range("H91").Value = -5800
$ws.Range("I91").Value = -4100
$ws.Range("J91").Value = -1800
$ws.Range("E94").Value = -764700
$ws.Range("F94").Value = -58300
$ws.Range("E102").Value = -232800
$ws.Range("F102").Value = 319300
